$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 110
$ws.Range("E12").Value = "110/140"
